# Apply the 2022-06-01 Fonds de solidarite Volet 1 data update.
# For each listed row, update the "nombre_aides" (column C) and
# "montant_total" (column E) values. "nombre_entreprises" (column D)
# is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 9;   C = 69575;  E = 191370177 },
    @{ Row = 13;  C = 37915;  E = 97638190 },
    @{ Row = 19;  C = 108918; E = 344661139 },
    @{ Row = 65;  C = 61056;  E = 334109625 },
    @{ Row = 157; C = 21203;  E = 77996841 },
    @{ Row = 164; C = 50577;  E = 168925084 },
    @{ Row = 168; C = 285014; E = 1210437825 },
    @{ Row = 169; C = 562606; E = 1284953022 },
    @{ Row = 170; C = 367399; E = 2845986686 },
    @{ Row = 171; C = 115164; E = 447107775 },
    @{ Row = 172; C = 21637;  E = 73049293 },
    @{ Row = 173; C = 54390;  E = 151889608 },
    @{ Row = 174; C = 357244; E = 1017908199 },
    @{ Row = 177; C = 96758;  E = 174746870 },
    @{ Row = 179; C = 235719; E = 812688286 },
    @{ Row = 180; C = 141486; E = 341034384 },
    @{ Row = 205; C = 11127;  E = 44175571 },
    @{ Row = 210; C = 6424;   E = 19297530 },
    @{ Row = 257; C = 182551; E = 1063828805 },
    @{ Row = 300; C = 15788;  E = 78280031 },
    @{ Row = 317; C = 103583; E = 303298714 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

$wb.Save()
